$wb = $excel.ActiveWorkbook

# Rename the first sheet: "Anthem - NH" -> "Anthem NH"
$wsAnthem = $wb.Worksheets.Item(1)
$wsAnthem.Name = "Anthem NH"

# Make "Anthem NH" the active/selected sheet (was "None", the 5th sheet)
# and move its selection to A23 (previously the whole used range A1:D21 was selected).
$wsAnthem.Activate()
$wsAnthem.Range("A23").Select()
